$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "surveys": lowercase / rename headers, style header row with a new
# black-font style, drop the extra S02 data row, fix selection
# ---------------------------------------------------------------------------
$wsSurveys = $wb.Worksheets.Item("surveys")
$wsSurveys.Range("B1").Value = "survey_id"
$wsSurveys.Range("E1").Value = "latitude"
$wsSurveys.Range("F1").Value = "longitude"
$wsSurveys.Range("A1:K1").Font.Color = 0x000000
$wsSurveys.Rows.Item(3).Delete()
$wsSurveys.Range("A1:K1").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "studies": lowercase the ID column header, move selection to A2
# ---------------------------------------------------------------------------
$wsStudies = $wb.Worksheets.Item("studies")
$wsStudies.Range("A1").Value = "study_id"
$wsStudies.Range("A2").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Notes": reword the explanatory note text on A3
# ---------------------------------------------------------------------------
$wsNotes = $wb.Worksheets.Item("Notes")
$wsNotes.Range("A3").Value = "Specific issue: survey_key contains IDs that are not present in surveys table"

# ---------------------------------------------------------------------------
# Sheet "counts": append the matching S02 row, move selection to A4
# ---------------------------------------------------------------------------
$wsCounts = $wb.Worksheets.Item("counts")
$wsCounts.Range("A3").Value = "S02"
$wsCounts.Range("B3").Value = "crt:1:A"
$wsCounts.Range("C3").Value = 1
$wsCounts.Range("D3").Value = 10
$wsCounts.Range("A4").Select() | Out-Null

# ---------------------------------------------------------------------------
# Make "Notes" the active tab/sheet (was "counts")
# ---------------------------------------------------------------------------
$wsNotes.Activate() | Out-Null
$wsNotes.Range("A4").Select() | Out-Null
